$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.55
$ws.Cells.Item(2, 9).Value = 3
$ws.Cells.Item(2, 12).Value = 3.75
$ws.Cells.Item(2, 26).Value = 26
$ws.Cells.Item(2, 27).Value = 26
$ws.Cells.Item(2, 54).Value = 301

# Row 4
$ws.Cells.Item(4, 13).Value = 1.11
$ws.Cells.Item(4, 14).Value = 6.5

# Row 8
$ws.Cells.Item(8, 7).Value = 1.4
$ws.Cells.Item(8, 8).Value = 4.33
$ws.Cells.Item(8, 9).Value = 8.5
$ws.Cells.Item(8, 10).Value = 1.91
$ws.Cells.Item(8, 11).Value = 2.38
$ws.Cells.Item(8, 12).Value = 7.5
$ws.Cells.Item(8, 13).Value = 1.04
$ws.Cells.Item(8, 14).Value = 13
$ws.Cells.Item(8, 17).Value = 1.83
$ws.Cells.Item(8, 18).Value = 2.03
$ws.Cells.Item(8, 21).Value = 2.1
$ws.Cells.Item(8, 22).Value = 1.67
$ws.Cells.Item(8, 24).Value = 6.5
$ws.Cells.Item(8, 25).Value = 9
$ws.Cells.Item(8, 26).Value = 9
$ws.Cells.Item(8, 27).Value = 12
$ws.Cells.Item(8, 28).Value = 29
$ws.Cells.Item(8, 30).Value = 8.5
$ws.Cells.Item(8, 31).Value = 21
$ws.Cells.Item(8, 32).Value = 67
$ws.Cells.Item(8, 33).Value = 351
$ws.Cells.Item(8, 34).Value = 19
$ws.Cells.Item(8, 35).Value = 41
$ws.Cells.Item(8, 36).Value = 23
$ws.Cells.Item(8, 37).Value = 101
$ws.Cells.Item(8, 40).Value = 3.25
$ws.Cells.Item(8, 41).Value = 7
$ws.Cells.Item(8, 42).Value = 21
$ws.Cells.Item(8, 43).Value = 21
$ws.Cells.Item(8, 47).Value = 9.5
$ws.Cells.Item(8, 48).Value = 67
$ws.Cells.Item(8, 49).Value = 8.5
$ws.Cells.Item(8, 50).Value = 41
$ws.Cells.Item(8, 52).Value = 151
$ws.Cells.Item(8, 53).Value = 201
$ws.Cells.Item(8, 54).Value = 351

# Row 10
$ws.Cells.Item(10, 7).Value = 1.73
$ws.Cells.Item(10, 8).Value = 3.3
$ws.Cells.Item(10, 9).Value = 5.25
$ws.Cells.Item(10, 10).Value = 2.5
$ws.Cells.Item(10, 12).Value = 4.75
$ws.Cells.Item(10, 13).Value = 1.08
$ws.Cells.Item(10, 14).Value = 8
$ws.Cells.Item(10, 19).Value = 1.4
$ws.Cells.Item(10, 20).Value = 2.75
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 1.73
$ws.Cells.Item(10, 24).Value = 7.5
$ws.Cells.Item(10, 26).Value = 13
$ws.Cells.Item(10, 36).Value = 17
$ws.Cells.Item(10, 40).Value = 3.75
$ws.Cells.Item(10, 41).Value = 10
$ws.Cells.Item(10, 42).Value = 21
$ws.Cells.Item(10, 43).Value = 34
$ws.Cells.Item(10, 45).Value = 151
$ws.Cells.Item(10, 46).Value = 2.75
$ws.Cells.Item(10, 47).Value = 8
$ws.Cells.Item(10, 48).Value = 51
$ws.Cells.Item(10, 49).Value = 6
$ws.Cells.Item(10, 50).Value = 23
$ws.Cells.Item(10, 51).Value = 34
$ws.Cells.Item(10, 52).Value = 81
$ws.Cells.Item(10, 53).Value = 101
$ws.Cells.Item(10, 54).Value = 251
